$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: copy the visual style (number format / font / fill / border /
# alignment) of one row's cells onto another row's cells, column by column,
# respecting the merged column groups used throughout the data table
# (A:B, C:G, H:K, L:M, N:O, P, Q).
# ---------------------------------------------------------------------------
function Copy-DataRowStyle($srcRow, $dstRow) {
    $groups = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
    foreach ($col in $groups) {
        $src = $ws.Range("$col$srcRow")
        $dst = $ws.Range("$col$dstRow")
        $dst.NumberFormat   = $src.NumberFormat
        $dst.Font.Name      = $src.Font.Name
        $dst.Font.Size      = $src.Font.Size
        $dst.Font.Color     = $src.Font.Color
        $dst.Interior.Color = $src.Interior.Color
        $dst.HorizontalAlignment = $src.HorizontalAlignment
        $dst.VerticalAlignment   = $src.VerticalAlignment
        $dst.WrapText       = $src.WrapText
        $dst.ShrinkToFit    = $src.ShrinkToFit
        $dst.Borders.LineStyle = $src.Borders.LineStyle
        $dst.Borders.Color     = $src.Borders.Color
    }
}

# ---------------------------------------------------------------------------
# Helper: write one product row's data + merges + row height.
# ---------------------------------------------------------------------------
function Set-ProductRow($row, $name, $balance, $orderLimit, $price, $sellPrice, $txCount, $height) {
    $ws.Range("C$row").Value = $name
    $ws.Range("H$row").Value = $balance
    # L/N/P hold numeric-looking text (e.g. "0", "45.00", "14.8500") that must
    # stay text (with its original trailing zeros) rather than be coerced to
    # a Number - prefix with an apostrophe to force text storage.
    $ws.Range("L$row").Value = "'" + $orderLimit
    $ws.Range("N$row").Value = "'" + $price
    $ws.Range("P$row").Value = "'" + $sellPrice
    $ws.Range("Q$row").Value = $txCount

    $ws.Range("A$row:B$row").Merge()
    $ws.Range("C$row:G$row").Merge()
    $ws.Range("H$row:K$row").Merge()
    $ws.Range("L$row:M$row").Merge()
    $ws.Range("N$row:O$row").Merge()

    $ws.Rows("$row`:$row").RowHeight = $height
}

# ---------------------------------------------------------------------------
# 1) Insert "CARVID 6.25MG 30TAB" as the new first data row (row 7), pushing
#    everything else down by one.
# ---------------------------------------------------------------------------
$ws.Rows("7:7").Insert()
Copy-DataRowStyle 8 7
Set-ProductRow 7 "CARVID 6.25MG 30TAB" "1:1" "0" "45.00" "14.8500" "0:1" 25.5
$ws.Range("A7").Value = 1

# ---------------------------------------------------------------------------
# 2) Insert "CYNCHOLINE 20CAPS." after CEPOREX (now row 8) -> new row 9.
# ---------------------------------------------------------------------------
$ws.Rows("9:9").Insert()
Copy-DataRowStyle 8 9
Set-ProductRow 9 "CYNCHOLINE 20CAPS." "0:1" "1" "20.00" "10.0000" "0:1" 25.5

# ---------------------------------------------------------------------------
# 3) Insert "URIVIN-N 10 EFF. SACHETS" before VOLTAREN (now row 12) -> new
#    row 12.
# ---------------------------------------------------------------------------
$ws.Rows("12:12").Insert()
Copy-DataRowStyle 11 12
Set-ProductRow 12 "URIVIN-N 10 EFF. SACHETS" "2:0" "1" "31.00" "31.0000" "1:0" 25.5

# ---------------------------------------------------------------------------
# 4) Renumber the "م" (row number) column for the whole table (1..8).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 8; $i++) {
    $ws.Range("A" + (7 + $i)).Value = $i + 1
}

# ---------------------------------------------------------------------------
# 5) Fix up the row-height alternation pattern (25.5 / 24.75) for every data
#    row, matching the final workbook exactly.
# ---------------------------------------------------------------------------
$heights = @{7=25.5; 8=24.75; 9=25.5; 10=24.75; 11=25.5; 12=25.5; 13=24.75; 14=25.5}
foreach ($r in $heights.Keys) {
    $ws.Rows("$r`:$r").RowHeight = $heights[$r]
}

# ---------------------------------------------------------------------------
# 6) Update the grand-total (sum of "سعر البيع" / sell price column) on the
#    summary row, which is now row 15.
# ---------------------------------------------------------------------------
$ws.Range("P15").Value = 142.4

# ---------------------------------------------------------------------------
# 7) Update the generated-on timestamp in the footer (now row 16).
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Saturday, 12 July, 2025 11:10 AM"
